$d = $word.ActiveDocument

# --- Step 1: Rewrite the text of the first paragraph ---------------------
# Old text: "test2: hitting on first read for some reason (valid bit should
#            be 0, so hit should be impossible)"
# New text: "Now writes and reads are hitting when the tag is the same, but
#            index is different"
$para1 = $d.Paragraphs(1).Range
$para1.Text = ""

$r = $d.Paragraphs(1).Range
$r.Collapse(1)
$r.InsertAfter("N")

$r = $d.Paragraphs(1).Range
$r.MoveEnd(1, -1)
$r.Collapse(0)
$r.InsertAfter("ow ")

$r = $d.Paragraphs(1).Range
$r.MoveEnd(1, -1)
$r.Collapse(0)
$r.InsertAfter("writes and reads are hitting when the tag is the same, but index is different")

# --- Step 2: Insert a brand-new paragraph right after paragraph 1 --------
# (after the _GoBack bookmark, before the "Write: Miss..." paragraph)
$p1 = $d.Paragraphs(1).Range
$p1.InsertParagraphAfter()

$newPara = $d.Paragraphs(2).Range
$newPara.Collapse(1)
$newPara.InsertAfter("test2: hitting on first read for some reason (valid bit should be 0, so hit should be impossible)")

$newParaRange = $d.Paragraphs(2).Range
$startPos = $newParaRange.Start
$firstRun = $d.Range($startPos, $startPos + 100)
$firstRun.HighlightColorIndex = 4

$r2 = $d.Paragraphs(2).Range
$r2.MoveEnd(1, -1)
$r2.Collapse(0)
$r2.InsertAfter(" ")

$newParaRange = $d.Paragraphs(2).Range
$spacePos = $newParaRange.Start + 100
$spaceRun = $d.Range($spacePos, $spacePos + 1)
$spaceRun.HighlightColorIndex = 4

$r3 = $d.Paragraphs(2).Range
$r3.MoveEnd(1, -1)
$r3.Collapse(0)
$r3.InsertAfter("Forgot to remove test line")

$newParaRange = $d.Paragraphs(2).Range
$lastPos = $newParaRange.Start + 101
$lastRun = $d.Range($lastPos, $lastPos + 27)
$lastRun.HighlightColorIndex = 0

Write-Output "Paragraph 1: [$($d.Paragraphs(1).Range.Text)]"
Write-Output "Paragraph 2: [$($d.Paragraphs(2).Range.Text)]"
